$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New transaction row (row 4): id=3, tanggal=2025-12-18, jenis=Income, kategori=a, jumlah=10000, catatan=(empty)
$ws.Range("A4").Value = 3
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "2025-12-18"
$ws.Range("C4").Value = "Income"
$ws.Range("D4").Value = "a"
$ws.Range("E4").Value = 10000
$ws.Range("F4").Value = ""
